# Adds new columns I (I0) and J (IF) to the sheet, mirroring column H's
# numeric series for rows 2..63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy formatting (bold, centered, bordered) from the
# existing "IP" header cell (H1), then set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..63 (column I and J share identical values,
# mirroring column H's series)
$values = @(9,9,8,9,8,8,8,7,8,8,8,9,7,8,8,8,8,8,7,10,8,7,8,9,9,8,9,8,8,8,8,8,8,9,9,10,8,8,8,8,8,8,8,7,8,8,7,7,8,8,8,7,8,8,9,5,8,6,6,8,5,4)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $v = $values[$i]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
